$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 177, shifting rows 177-234
# down to 178-235 (keeps every other column's values attached to their
# original row, matching the diff's row-by-row value shuffle).
$ws.Rows.Item(177).Insert()

$newRow = 177
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"
$ws.Cells.Item($newRow, 4).Value = 44559
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100114014
$ws.Cells.Item($newRow, 7).Value = "Betarraga"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 5000
$ws.Cells.Item($newRow, 11).Value = 500
$ws.Cells.Item($newRow, 12).Value = 500
$ws.Cells.Item($newRow, 13).Value = 500
$ws.Cells.Item($newRow, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 100
$ws.Cells.Item($newRow, 17).Value = 5
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
